$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2126.3333
$ws.Range("J17").Value = 2126.3333
$ws.Range("L17").Value = 6378.999899999999
$ws.Range("N17").Value = -6714.999899999999
$ws.Range("H92").Value = 648.75
$ws.Range("I92").Value = 698.3333
$ws.Range("K92").Value = 698.3333
$ws.Range("M92").Value = 549.6667
$ws.Range("H104").Value = 1200
$ws.Range("I104").Value = 1200
$ws.Range("K104").Value = 3600
$ws.Range("M104").Value = -1853
$ws.Range("H112").Value = 3444.4443
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 3444.4443
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 10333.3329
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -12549.3329
$ws.Range("H121").Value = 2512.4285
$ws.Range("J121").Value = 2512.4285
$ws.Range("L121").Value = 7537.2855
$ws.Range("N121").Value = -11031.2855
$ws.Range("H137").Value = 7278.294
$ws.Range("I137").Value = 7324
$ws.Range("K137").Value = 21972
$ws.Range("M137").Value = -19422
$ws.Range("H138").Value = 6006.364
$ws.Range("I138").Value = 3683.5557
$ws.Range("J138").Value = 6373.123
$ws.Range("K138").Value = 11050.6671
$ws.Range("L138").Value = 19119.369
$ws.Range("M138").Value = -5910.667099999999
$ws.Range("N138").Value = -29399.369

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3593.8354
$ws.Range("I32").Value = 3349.3506
$ws.Range("K32").Value = 3349.3506
$ws.Range("M32").Value = -3062.3506
$ws.Range("H37").Value = 58950
$ws.Range("J37").Value = 58950
$ws.Range("L37").Value = 58950
$ws.Range("N37").Value = -59496
$ws.Range("H55").Value = 29276.5
$ws.Range("J55").Value = 29276.5
$ws.Range("L55").Value = 29276.5
$ws.Range("N55").Value = -29906.5
$ws.Range("H61").Value = 4495.5
$ws.Range("I61").Value = 4207.7
$ws.Range("J61").Value = 5934.5
$ws.Range("K61").Value = 4207.7
$ws.Range("L61").Value = 5934.5
$ws.Range("M61").Value = -3995.7
$ws.Range("N61").Value = -6358.5
$ws.Range("H74").Value = 2054.1428
$ws.Range("I74").Value = 1746.3334
$ws.Range("J74").Value = 2608.2
$ws.Range("K74").Value = 1746.3334
$ws.Range("L74").Value = 2608.2
$ws.Range("M74").Value = -872.3334
$ws.Range("N74").Value = -4356.2
$ws.Range("H77").Value = 2054.1428
$ws.Range("I77").Value = 1746.3334
$ws.Range("J77").Value = 2608.2
$ws.Range("K77").Value = 8731.666999999999
$ws.Range("L77").Value = 13041
$ws.Range("M77").Value = -4363.666999999999
$ws.Range("N77").Value = -21777
$ws.Range("H97").Value = 954.5238000000001
$ws.Range("I97").Value = 954.5238000000001
$ws.Range("K97").Value = 954.5238000000001
$ws.Range("M97").Value = -458.5238000000001
$ws.Range("H135").Value = 80809.336
$ws.Range("J135").Value = 80809.336
$ws.Range("L135").Value = 80809.336
$ws.Range("N135").Value = -90949.336
$ws.Range("H136").Value = 4495.5
$ws.Range("I136").Value = 4207.7
$ws.Range("J136").Value = 5934.5
$ws.Range("K136").Value = 12623.1
$ws.Range("L136").Value = 17803.5
$ws.Range("M136").Value = -10073.1
$ws.Range("N136").Value = -22903.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2883.9
$ws.Range("I20").Value = 2714
$ws.Range("K20").Value = 2714
$ws.Range("M20").Value = -2467
$ws.Range("H86").Value = 1704097.1
$ws.Range("I86").Value = 2128874.5
$ws.Range("K86").Value = 2128874.5
$ws.Range("M86").Value = -2127751.5
$ws.Range("H89").Value = 1704097.1
$ws.Range("I89").Value = 2128874.5
$ws.Range("K89").Value = 10644372.5
$ws.Range("M89").Value = -10638756.5
$ws.Range("H94").Value = 1414.9048
$ws.Range("I94").Value = 1534.6111
$ws.Range("K94").Value = 1534.6111
$ws.Range("M94").Value = -1083.6111
$ws.Range("H99").Value = 2124.4666
$ws.Range("I99").Value = 1875.4615
$ws.Range("J99").Value = 3743
$ws.Range("K99").Value = 1875.4615
$ws.Range("L99").Value = 3743
$ws.Range("M99").Value = -377.4614999999999
$ws.Range("N99").Value = -6739
$ws.Range("H134").Value = 48671.957
$ws.Range("I134").Value = 6906.45
$ws.Range("K134").Value = 20719.35
$ws.Range("M134").Value = -18184.35

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 20252000
$ws.Range("I4").Value = 501500.5
$ws.Range("K4").Value = 501500.5
$ws.Range("M4").Value = -501388.5
$ws.Range("H16").Value = 8193.352999999999
$ws.Range("I16").Value = 3709.3
$ws.Range("K16").Value = 3709.3
$ws.Range("M16").Value = -3422.3
$ws.Range("H22").Value = 409.9
$ws.Range("I22").Value = 344.33334
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 344.33334
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = 5.666659999999979
$ws.Range("N22").Value = -1700
$ws.Range("H62").Value = 3765.8572
$ws.Range("I62").Value = 2691
$ws.Range("J62").Value = 6453
$ws.Range("K62").Value = 2691
$ws.Range("L62").Value = 6453
$ws.Range("M62").Value = -2067
$ws.Range("N62").Value = -7701
$ws.Range("H65").Value = 3765.8572
$ws.Range("I65").Value = 2691
$ws.Range("J65").Value = 6453
$ws.Range("K65").Value = 13455
$ws.Range("L65").Value = 32265
$ws.Range("M65").Value = -10335
$ws.Range("N65").Value = -38505
$ws.Range("H107").Value = 508.16
$ws.Range("I107").Value = 443.8421
$ws.Range("J107").Value = 711.8333
$ws.Range("K107").Value = 443.8421
$ws.Range("L107").Value = 711.8333
$ws.Range("M107").Value = 1476.1579
$ws.Range("N107").Value = -4551.8333
$ws.Range("H113").Value = 8193.352999999999
$ws.Range("I113").Value = 3709.3
$ws.Range("K113").Value = 3709.3
$ws.Range("M113").Value = -1539.3
$ws.Range("H122").Value = 4258.7617
$ws.Range("I122").Value = 3536.0833
$ws.Range("K122").Value = 10608.2499
$ws.Range("M122").Value = -8158.249899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 1979.4
$ws.Range("J124").Value = 2299
$ws.Range("L124").Value = 6897
$ws.Range("N124").Value = -16717
$ws.Range("H125").Value = 4915
$ws.Range("I125").Value = 4030
$ws.Range("J125").Value = 5800
$ws.Range("K125").Value = 12090
$ws.Range("L125").Value = 17400
$ws.Range("M125").Value = -7170
$ws.Range("N125").Value = -27240
$ws.Range("H131").Value = 5477.35
$ws.Range("I131").Value = 1677.7778
$ws.Range("J131").Value = 8586.091
$ws.Range("K131").Value = 5033.3334
$ws.Range("L131").Value = 25758.273
$ws.Range("M131").Value = 6.666599999999562
$ws.Range("N131").Value = -35838.273
$ws.Range("H140").Value = 3917.44
$ws.Range("I140").Value = 2260.9412
$ws.Range("J140").Value = 7437.5
$ws.Range("K140").Value = 6782.823600000001
$ws.Range("L140").Value = 22312.5
$ws.Range("M140").Value = -1602.823600000001
$ws.Range("N140").Value = -32672.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 17149996
$ws.Range("I33").Value = 29989
$ws.Range("K33").Value = 29989
$ws.Range("M33").Value = -29737
$ws.Range("H113").Value = 628167.1
$ws.Range("I113").Value = 1430184.6
$ws.Range("K113").Value = 1430184.6
$ws.Range("M113").Value = -1428014.6
$ws.Range("H123").Value = 52496
$ws.Range("J123").Value = 52496
$ws.Range("L123").Value = 52496
$ws.Range("N123").Value = -57396

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1239.0667
$ws.Range("I22").Value = 1100
$ws.Range("K22").Value = 1100
$ws.Range("M22").Value = -805
$ws.Range("H27").Value = 1239.0667
$ws.Range("I27").Value = 1100
$ws.Range("K27").Value = 1100
$ws.Range("M27").Value = -993
$ws.Range("H55").Value = 2797.3333
$ws.Range("I55").Value = 213.33333
$ws.Range("K55").Value = 213.33333
$ws.Range("M55").Value = -40.33332999999999
$ws.Range("H61").Value = 3289.5217
$ws.Range("I61").Value = 2912.2354
$ws.Range("J61").Value = 4358.5
$ws.Range("K61").Value = 2912.2354
$ws.Range("L61").Value = 4358.5
$ws.Range("M61").Value = -2710.2354
$ws.Range("N61").Value = -4762.5
$ws.Range("H93").Value = 62502830
$ws.Range("I93").Value = 100001950
$ws.Range("J93").Value = 4300.3335
$ws.Range("K93").Value = 100001950
$ws.Range("L93").Value = 4300.3335
$ws.Range("M93").Value = -100000702
$ws.Range("N93").Value = -6796.3335
$ws.Range("H113").Value = 3289.5217
$ws.Range("I113").Value = 2912.2354
$ws.Range("J113").Value = 4358.5
$ws.Range("K113").Value = 2912.2354
$ws.Range("L113").Value = 4358.5
$ws.Range("M113").Value = -742.2354
$ws.Range("N113").Value = -8698.5
$ws.Range("H122").Value = 3714.6072
$ws.Range("I122").Value = 2476.7144
$ws.Range("K122").Value = 7430.1432
$ws.Range("M122").Value = -4980.1432
$ws.Range("H132").Value = 8919.727999999999
$ws.Range("I132").Value = 7525.7
$ws.Range("J132").Value = 10081.417
$ws.Range("K132").Value = 22577.1
$ws.Range("L132").Value = 30244.251
$ws.Range("M132").Value = -20047.1
$ws.Range("N132").Value = -35304.251
$ws.Range("H136").Value = 269392.4
$ws.Range("I136").Value = 504320.16
$ws.Range("K136").Value = 1512960.48
$ws.Range("M136").Value = -1510410.48

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1876
$ws.Range("J126").Value = 1869
$ws.Range("L126").Value = 5607
$ws.Range("N126").Value = -10547
$ws.Range("H132").Value = 31267.324
$ws.Range("I132").Value = 1942.2778
$ws.Range("J132").Value = 59048.95
$ws.Range("K132").Value = 5826.8334
$ws.Range("L132").Value = 177146.85
$ws.Range("M132").Value = -3296.8334
$ws.Range("N132").Value = -182206.85
$ws.Range("H136").Value = 10832624
$ws.Range("I136").Value = 12838486
$ws.Range("K136").Value = 38515458
$ws.Range("M136").Value = -38512908
